$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 27177.25
$ws.Range("I57").Value = 19354.5
$ws.Range("J57").Value = 35000
$ws.Range("K57").Value = 58063.5
$ws.Range("L57").Value = 105000
$ws.Range("M57").Value = -57564.5
$ws.Range("N57").Value = -105998

$ws.Range("H82").Value = 6295.5713
$ws.Range("I82").Value = 1605.8572
$ws.Range("J82").Value = 10985.286
$ws.Range("K82").Value = 4817.571599999999
$ws.Range("L82").Value = 32955.858
$ws.Range("M82").Value = -4411.571599999999

$ws.Range("H85").Value = 6295.5713
$ws.Range("I85").Value = 1605.8572
$ws.Range("J85").Value = 10985.286
$ws.Range("K85").Value = 4817.571599999999
$ws.Range("L85").Value = 32955.858
$ws.Range("M85").Value = -3413.571599999999

$ws.Range("H113").Value = 7860.4
$ws.Range("I113").Value = 4457.4287
$ws.Range("J113").Value = 9692.77
$ws.Range("K113").Value = 4457.4287
$ws.Range("L113").Value = 9692.77
$ws.Range("M113").Value = -1203.4287
$ws.Range("N113").Value = -16200.77

$ws.Range("H132").Value = 3840.7097
$ws.Range("I132").Value = 2736.6667
$ws.Range("J132").Value = 6159.2
$ws.Range("K132").Value = 8210.000100000001
$ws.Range("L132").Value = 18477.6
$ws.Range("M132").Value = -5680.000100000001
$ws.Range("N132").Value = -23537.6

$ws.Range("H138").Value = 3065.4084
$ws.Range("I138").Value = 1585.875
$ws.Range("J138").Value = 3495.818
$ws.Range("K138").Value = 4757.625
$ws.Range("L138").Value = 10487.454
$ws.Range("M138").Value = 382.375
$ws.Range("N138").Value = -20767.454

$ws.Range("H141").Value = 73914.36
$ws.Range("I141").Value = 92991
$ws.Range("J141").Value = 3966.6667
$ws.Range("K141").Value = 278973
$ws.Range("L141").Value = 11900.0001
$ws.Range("M141").Value = -273793
$ws.Range("N141").Value = -22260.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 904.1724
$ws.Range("I2").Value = 922.4211
$ws.Range("J2").Value = 869.5
$ws.Range("K2").Value = 922.4211
$ws.Range("L2").Value = 869.5
$ws.Range("M2").Value = -809.4211
$ws.Range("N2").Value = -1095.5

$ws.Range("H45").Value = 1147.7142
$ws.Range("I45").Value = 1000
$ws.Range("J45").Value = 1172.3334
$ws.Range("K45").Value = 1000
$ws.Range("L45").Value = 1172.3334
$ws.Range("M45").Value = -623
$ws.Range("N45").Value = -1926.3334

$ws.Range("H57").Value = 3000
$ws.Range("I57").Value = 3000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 3000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -2516

$ws.Range("H61").Value = 2481.5557
$ws.Range("I61").Value = 1978.5714
$ws.Range("J61").Value = 2801.6365
$ws.Range("K61").Value = 1978.5714
$ws.Range("L61").Value = 2801.6365
$ws.Range("M61").Value = -1766.5714
$ws.Range("N61").Value = -3225.6365

$ws.Range("H74").Value = 3469.1667
$ws.Range("I74").Value = 3612.697
$ws.Range("J74").Value = 2942.889
$ws.Range("K74").Value = 3612.697
$ws.Range("L74").Value = 2942.889
$ws.Range("M74").Value = -2738.697

$ws.Range("H77").Value = 3469.1667
$ws.Range("I77").Value = 3612.697
$ws.Range("J77").Value = 2942.889
$ws.Range("K77").Value = 18063.485
$ws.Range("L77").Value = 14714.445
$ws.Range("M77").Value = -13695.485

$ws.Range("H112").Value = 26238.096
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 26238.096
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 26238.096
$ws.Range("N112").Value = -29192.096

$ws.Range("H116").Value = 904.1724
$ws.Range("I116").Value = 922.4211
$ws.Range("J116").Value = 869.5
$ws.Range("K116").Value = 922.4211
$ws.Range("L116").Value = 869.5
$ws.Range("M116").Value = 1371.5789
$ws.Range("N116").Value = -5457.5

$ws.Range("H122").Value = 1872.6
$ws.Range("I122").Value = 1391.1904
$ws.Range("J122").Value = 4400
$ws.Range("K122").Value = 4173.5712
$ws.Range("L122").Value = 13200
$ws.Range("M122").Value = -1723.5712
$ws.Range("N122").Value = -18100

$ws.Range("H124").Value = 29419
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 29419
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 29419
$ws.Range("N124").Value = -39239

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 2480.4666
$ws.Range("I132").Value = 1914.2354
$ws.Range("J132").Value = 4230.636
$ws.Range("K132").Value = 5742.706200000001
$ws.Range("L132").Value = 12691.908
$ws.Range("M132").Value = -3212.706200000001
$ws.Range("N132").Value = -17751.908

$ws.Range("H136").Value = 2481.5557
$ws.Range("I136").Value = 1978.5714
$ws.Range("J136").Value = 2801.6365
$ws.Range("K136").Value = 5935.7142
$ws.Range("L136").Value = 8404.9095
$ws.Range("M136").Value = -3385.7142
$ws.Range("N136").Value = -13504.9095

$ws.Range("H137").Value = 39774.285
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 39774.285
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 39774.285
$ws.Range("N137").Value = -49974.285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 904.1724
$ws.Range("I3").Value = 922.4211
$ws.Range("J3").Value = 869.5
$ws.Range("K3").Value = 922.4211
$ws.Range("L3").Value = 869.5
$ws.Range("M3").Value = -808.4211
$ws.Range("N3").Value = -1097.5

$ws.Range("H59").Value = 81998
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 81998
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 81998
$ws.Range("N59").Value = -83692

$ws.Range("H99").Value = 4848.2144
$ws.Range("I99").Value = 1310
$ws.Range("J99").Value = 5813.1816
$ws.Range("K99").Value = 1310
$ws.Range("L99").Value = 5813.1816
$ws.Range("M99").Value = 188
$ws.Range("N99").Value = -8809.1816

$ws.Range("H112").Value = 30417.5
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 30417.5
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 30417.5
$ws.Range("N112").Value = -33371.5

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("M128").ClearContents()

$ws.Range("H137").Value = 36205
$ws.Range("I137").Value = 25000
$ws.Range("J137").Value = 39940
$ws.Range("K137").Value = 25000
$ws.Range("L137").Value = 39940
$ws.Range("M137").Value = -19900
$ws.Range("N137").Value = -50140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 627
$ws.Range("I107").Value = 269.36365
$ws.Range("J107").Value = 1282.6666
$ws.Range("K107").Value = 269.36365
$ws.Range("L107").Value = 1282.6666
$ws.Range("M107").Value = 1650.63635
$ws.Range("N107").Value = -5122.6666

$ws.Range("H134").Value = 2276.3667
$ws.Range("I134").Value = 1190.8
$ws.Range("J134").Value = 3361.9333
$ws.Range("K134").Value = 3572.4
$ws.Range("L134").Value = 10085.7999
$ws.Range("M134").Value = -1037.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1611.0358
$ws.Range("I132").Value = 795.6
$ws.Range("J132").Value = 2551.923
$ws.Range("K132").Value = 7160.400000000001
$ws.Range("L132").Value = 22967.307
$ws.Range("M132").Value = -4630.400000000001
$ws.Range("N132").Value = -28027.307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 39660
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 39660
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 39660
$ws.Range("N100").Value = -41824

$ws.Range("H111").Value = 24500
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 24500
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 24500
$ws.Range("N111").Value = -30634

$ws.Range("H122").Value = 3238.0527
$ws.Range("I122").Value = 1536.7858
$ws.Range("J122").Value = 8001.6
$ws.Range("K122").Value = 4610.357400000001
$ws.Range("L122").Value = 24004.8
$ws.Range("M122").Value = -2160.357400000001

$ws.Range("H126").Value = 2909.09
$ws.Range("I126").Value = 2918.2727
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 8754.8181
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -6284.8181
$ws.Range("N126").Value = -10940

$ws.Range("H137").Value = 42696.668
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 42696.668
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 42696.668
$ws.Range("N137").Value = -52896.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 35739.6
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 35739.6
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 35739.6
$ws.Range("N110").Value = -43919.6

$ws.Range("H122").Value = 5815.857
$ws.Range("I122").Value = 3615.2856
$ws.Range("J122").Value = 8016.4287
$ws.Range("K122").Value = 10845.8568
$ws.Range("L122").Value = 24049.2861
$ws.Range("M122").Value = -8395.856800000001
$ws.Range("N122").Value = -28949.2861

$ws.Range("H132").Value = 3591.9824
$ws.Range("I132").Value = 1114.1025
$ws.Range("J132").Value = 8960.723
$ws.Range("K132").Value = 3342.3075
$ws.Range("L132").Value = 26882.169
$ws.Range("M132").Value = -812.3074999999999
$ws.Range("N132").Value = -31942.169

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H113").Value = 7569.5
$ws.Range("I113").Value = 9370.272000000001
$ws.Range("J113").Value = 966.6667
$ws.Range("K113").Value = 28110.816
$ws.Range("L113").Value = 2900.0001
$ws.Range("M113").Value = -25940.816

$ws.Range("H122").Value = 4671.2
$ws.Range("I122").Value = 2758.5
$ws.Range("J122").Value = 6857.143
$ws.Range("K122").Value = 8275.5
$ws.Range("L122").Value = 20571.429
$ws.Range("M122").Value = -5825.5
$ws.Range("N122").Value = -25471.429

$ws.Range("H123").Value = 37888.625
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 37888.625
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 37888.625
$ws.Range("N123").Value = -47688.625

$ws.Range("H128").Value = 41795
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 41795
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 41795
$ws.Range("N128").Value = -51755

$ws.Range("H132").Value = 20842438
$ws.Range("I132").Value = 15285.429
$ws.Range("J132").Value = 37041336
$ws.Range("K132").Value = 45856.287
$ws.Range("L132").Value = 111124008
$ws.Range("M132").Value = -43326.287
$ws.Range("N132").Value = -111129068
